$d = $word.ActiveDocument

# --- Change 1: append an emailing contact after "...questions they may have." ---
$d.Content.Find.Execute(
    " have.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " have, by emailing Melissa Kline at Melissa.e.kline@gmail.com ",
    2) | Out-Null

# --- Change 2: move the (hidden) "_GoBack" bookmark from its old spot -----
# (between "...about what your" and " children think...") to a new spot
# that splits "animated videos" into "animated video" | "s of a cartoon...".

# First, collapse the old bookmark location back into plain running text
# (re-writing the same text over the bookmark removes it and re-merges the
# two runs it used to separate).
$d.Content.Find.Execute(
    "your children think",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "your children think",
    2) | Out-Null

# Now add the "_GoBack" bookmark at its new location; Word only ever keeps
# one bookmark with a given name, so this implicitly finishes the move.
$text = $d.Content.Text
$marker = "During the study, your child will see animated video"
$pos = $text.IndexOf($marker)
$splitAt = $pos + $marker.Length
$bookmarkRange = $d.Range($splitAt, $splitAt)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null
